# Weekly data refresh: a new price record for "Brócoli" (Terminal
# Hortofrutícola Agro Chillán) is published at the top of the historical
# block, pushing every existing record in that block down by one row.
#
# The data block for this report lives at rows 287-343 (header is row 1,
# data starts at row 2). Inserting a whole row at 287 shifts the existing
# 287-343 records down to 288-344 (carrying their values/number formats
# with them, exactly like Excel's native "Insert Copied/Sheet Rows"), and
# the new weekly record is then written into the now-empty row 287.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 287:343 down to 288:344, leaving a fresh (empty) row 287.
$ws.Rows.Item(287).Insert()

# Populate the new record in row 287.
$ws.Cells.Item(287, 1).Value = 7
$ws.Cells.Item(287, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(287, 3).Value = "Ñuble"
$ws.Cells.Item(287, 4).Value = 44889
$ws.Cells.Item(287, 5).Value = 16
$ws.Cells.Item(287, 6).Value = 100112023
$ws.Cells.Item(287, 7).Value = "Brócoli"
$ws.Cells.Item(287, 8).Value = "Sin especificar"
$ws.Cells.Item(287, 9).Value = "Segunda"
$ws.Cells.Item(287, 10).Value = 200
$ws.Cells.Item(287, 11).Value = 600
$ws.Cells.Item(287, 12).Value = 600
$ws.Cells.Item(287, 13).Value = 600
$ws.Cells.Item(287, 14).Value = "`$/unidad"
$ws.Cells.Item(287, 15).Value = "Región del Maule"
$ws.Cells.Item(287, 16).Value = 600
$ws.Cells.Item(287, 17).Value = 1
$ws.Cells.Item(287, 18).Value = "Hortaliza"
